# Spain Primera Liga workbook update
# Commit: Atualizacao de bases das ligas, do dia: 2024-01-29 as 07-55
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: Row content rotations (swap match rows B:AC, keep id column A fixed) ---
# Rows [228, 229]
$v228 = $ws.Range("B228:AC228").Value2
$v229 = $ws.Range("B229:AC229").Value2
$ws.Range("B228:AC228").Value2 = $v229
$ws.Range("B229:AC229").Value2 = $v228

# Rows [424, 425]
$v424 = $ws.Range("B424:AC424").Value2
$v425 = $ws.Range("B425:AC425").Value2
$ws.Range("B424:AC424").Value2 = $v425
$ws.Range("B425:AC425").Value2 = $v424

# Rows [580, 581]
$v580 = $ws.Range("B580:AC580").Value2
$v581 = $ws.Range("B581:AC581").Value2
$ws.Range("B580:AC580").Value2 = $v581
$ws.Range("B581:AC581").Value2 = $v580

# Rows [597, 598]
$v597 = $ws.Range("B597:AC597").Value2
$v598 = $ws.Range("B598:AC598").Value2
$ws.Range("B597:AC597").Value2 = $v598
$ws.Range("B598:AC598").Value2 = $v597

# Rows [604, 605, 606]
$v604 = $ws.Range("B604:AC604").Value2
$v605 = $ws.Range("B605:AC605").Value2
$v606 = $ws.Range("B606:AC606").Value2
$ws.Range("B604:AC604").Value2 = $v606
$ws.Range("B605:AC605").Value2 = $v604
$ws.Range("B606:AC606").Value2 = $v605

# Rows [755, 756]
$v755 = $ws.Range("B755:AC755").Value2
$v756 = $ws.Range("B756:AC756").Value2
$ws.Range("B755:AC755").Value2 = $v756
$ws.Range("B756:AC756").Value2 = $v755

# Rows [912, 913]
$v912 = $ws.Range("B912:AC912").Value2
$v913 = $ws.Range("B913:AC913").Value2
$ws.Range("B912:AC912").Value2 = $v913
$ws.Range("B913:AC913").Value2 = $v912

# Rows [931, 932]
$v931 = $ws.Range("B931:AC931").Value2
$v932 = $ws.Range("B932:AC932").Value2
$ws.Range("B931:AC931").Value2 = $v932
$ws.Range("B932:AC932").Value2 = $v931

# Rows [935, 936]
$v935 = $ws.Range("B935:AC935").Value2
$v936 = $ws.Range("B936:AC936").Value2
$ws.Range("B935:AC935").Value2 = $v936
$ws.Range("B936:AC936").Value2 = $v935

# Rows [957, 958]
$v957 = $ws.Range("B957:AC957").Value2
$v958 = $ws.Range("B958:AC958").Value2
$ws.Range("B957:AC957").Value2 = $v958
$ws.Range("B958:AC958").Value2 = $v957

# Rows [972, 973, 974, 975]
$v972 = $ws.Range("B972:AC972").Value2
$v973 = $ws.Range("B973:AC973").Value2
$v974 = $ws.Range("B974:AC974").Value2
$v975 = $ws.Range("B975:AC975").Value2
$ws.Range("B972:AC972").Value2 = $v973
$ws.Range("B973:AC973").Value2 = $v975
$ws.Range("B974:AC974").Value2 = $v972
$ws.Range("B975:AC975").Value2 = $v974

# Rows [1050, 1051]
$v1050 = $ws.Range("B1050:AC1050").Value2
$v1051 = $ws.Range("B1051:AC1051").Value2
$ws.Range("B1050:AC1050").Value2 = $v1051
$ws.Range("B1051:AC1051").Value2 = $v1050

# Rows [1143, 1144]
$v1143 = $ws.Range("B1143:AC1143").Value2
$v1144 = $ws.Range("B1144:AC1144").Value2
$ws.Range("B1143:AC1143").Value2 = $v1144
$ws.Range("B1144:AC1144").Value2 = $v1143

# Rows [1156, 1157]
$v1156 = $ws.Range("B1156:AC1156").Value2
$v1157 = $ws.Range("B1157:AC1157").Value2
$ws.Range("B1156:AC1156").Value2 = $v1157
$ws.Range("B1157:AC1157").Value2 = $v1156

# Rows [1159, 1160]
$v1159 = $ws.Range("B1159:AC1159").Value2
$v1160 = $ws.Range("B1160:AC1160").Value2
$ws.Range("B1159:AC1159").Value2 = $v1160
$ws.Range("B1160:AC1160").Value2 = $v1159

# Rows [1161, 1162]
$v1161 = $ws.Range("B1161:AC1161").Value2
$v1162 = $ws.Range("B1162:AC1162").Value2
$ws.Range("B1161:AC1161").Value2 = $v1162
$ws.Range("B1162:AC1162").Value2 = $v1161

# --- Part 2: Individual odds-value edits (rows 1202-1204) ---
$ws.Range("N1202").Value2 = 1.833
$ws.Range("P1202").Value2 = 4.75
$ws.Range("R1202").Value2 = 1.85
$ws.Range("S1202").Value2 = 2.08
$ws.Range("U1202").Value2 = 2.05
$ws.Range("V1202").Value2 = 1.85
$ws.Range("O1203").Value2 = 5
$ws.Range("P1203").Value2 = 9
$ws.Range("U1203").Value2 = 2.06
$ws.Range("V1203").Value2 = 1.84
$ws.Range("R1204").Value2 = 2.05
$ws.Range("S1204").Value2 = 1.85
$ws.Range("U1204").Value2 = 1.92
$ws.Range("V1204").Value2 = 1.98

# --- Part 3: Append new row 1215 (future/unplayed match) ---
$ws.Range("A1214").Copy($ws.Range("A1215"))
$ws.Range("E1214").Copy($ws.Range("E1215"))

$ws.Cells.Item(1215,1).Value2 = 1213
$ws.Cells.Item(1215,2).Value2 = 6809314
$ws.Cells.Item(1215,3).Value2 = "Spain Primera Liga"
$ws.Cells.Item(1215,4).Value2 = "Spain Primera Liga"
$ws.Cells.Item(1215,5).Value2 = 45327.70833333334
$ws.Cells.Item(1215,6).Value2 = "Rayo Vallecano"
$ws.Cells.Item(1215,7).Value2 = "Sevilla"
$ws.Cells.Item(1215,11).Value2 = 2.375
$ws.Cells.Item(1215,12).Value2 = 3.25
$ws.Cells.Item(1215,13).Value2 = 3
$ws.Cells.Item(1215,14).Value2 = 2.2
$ws.Cells.Item(1215,15).Value2 = 3.3
$ws.Cells.Item(1215,16).Value2 = 3.3
$ws.Cells.Item(1215,17).Value2 = -0.25
$ws.Cells.Item(1215,18).Value2 = 1.92
$ws.Cells.Item(1215,19).Value2 = 1.98
$ws.Cells.Item(1215,20).Value2 = 2.25
$ws.Cells.Item(1215,21).Value2 = 1.99
$ws.Cells.Item(1215,22).Value2 = 1.91
$ws.Cells.Item(1215,23).Value2 = 0
$ws.Cells.Item(1215,24).Value2 = 0
$ws.Cells.Item(1215,25).Value2 = 0
$ws.Cells.Item(1215,26).Value2 = 0
$ws.Cells.Item(1215,27).Value2 = 0
